$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- Insert the new "GRANITO" column at F (col 6) ---
# This shifts old F..N (GRIFERIAS..SAL SOLUBLE) to G..O, and Insert()
# naturally copies formatting from the left neighbour column so the
# header (s=5), data (s=1) and totals-row (s=3) styles all come along.
$ws.Columns.Item(6).Insert()

# --- Append three new trailing columns: P, Q, R ---
# Inserting right after the current last column (O = 15) still clones
# the neighbouring format, the same way a mid-sheet insert does.
$ws.Columns.Item(16).Insert()
$ws.Columns.Item(17).Insert()
$ws.Columns.Item(18).Insert()

# --- Column widths (stored OOXML "width", not the ColumnWidth char-width) ---
# ColumnWidth and the saved <col width="..."> differ by the standard
# ~0.83 padding offset used by this workbook, so subtract it back out.
$ws.Columns.Item(6).ColumnWidth = 12.17   # F  -> stored width 13
$ws.Columns.Item(16).ColumnWidth = 19.17  # P  -> stored width 20
$ws.Columns.Item(17).ColumnWidth = 16.17  # Q  -> stored width 17
$ws.Columns.Item(18).ColumnWidth = 15.17  # R  -> stored width 16

# --- Header row ---
$ws.Range("F1").Value = "GRANITO"
$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

# --- Data rows 2-29 ---
# GRANITO is an entirely new column with no prior sales, so it's all 0.
$ws.Range("F2:F29").Value = 0
$ws.Range("P2:P29").Value = 0
$ws.Range("Q2:Q29").Value = 0
$ws.Range("R2:R29").Value = 0

# A couple of rows have real NO RESURTIBLES figures.
$ws.Range("P17").Value = 129.66
$ws.Range("P29").Value = 41.6

# --- Totals row 30 ("x de 28" counters) ---
$ws.Range("F30").Value = "0 de 28"
$ws.Range("P30").Value = "2 de 28"
$ws.Range("Q30").Value = "0 de 28"
$ws.Range("R30").Value = "0 de 28"
